$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the Mobile No (D2:D11) values - they will now be populated via Twilio API
$ws.Range("D2:D11").ClearContents()

# Format column D: right-aligned
$rng = $ws.Range("D2:D11")
$rng.HorizontalAlignment = -4152  # xlRight

# Set column D width (best fit-like width)
$ws.Columns("D").ColumnWidth = 13.42578125

# Update selection
$ws.Range("D2:D11").Select()
